# Add 9 days of missing historical data (2019-11-18 .. 2019-11-28) right
# after the existing 2019-11-15 row, shifting all subsequent rows down by 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows starting at row 749 (pushes old row 749.. down to 758..)
$ws.Rows("749:757").Insert()

# Force columns B (date) and C (id) to text so strings like "2019-11-18"
# and "0187" are not auto-converted to a date serial / number.
$ws.Range("B749:B757").NumberFormat = "@"
$ws.Range("C749:C757").NumberFormat = "@"

$newRows = @(
    @(1574035200, "2019-11-18", 0.14,  0.145, 0.14,  0.145, 6100),
    @(1574121600, "2019-11-19", 0.14,  0.145, 0.14,  0.14,  536900),
    @(1574208000, "2019-11-20", 0.14,  0.145, 0.14,  0.145, 702700),
    @(1574294400, "2019-11-21", 0.145, 0.145, 0.14,  0.145, 540800),
    @(1574380800, "2019-11-22", 0.145, 0.145, 0.135, 0.14,  1571500),
    @(1574640000, "2019-11-25", 0.14,  0.14,  0.135, 0.135, 378100),
    @(1574726400, "2019-11-26", 0.135, 0.14,  0.135, 0.14,  81000),
    @(1574812800, "2019-11-27", 0.135, 0.135, 0.13,  0.135, 1651700),
    @(1574899200, "2019-11-28", 0.135, 0.135, 0.135, 0.135, 135600)
)

$r = 749
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "0187"
    $ws.Cells.Item($r, 4).Value = "BCMALL"
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $ws.Cells.Item($r, 9).Value = $row[6]
    $r = $r + 1
}
